$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 6820800
$ws.Range("E8").Value = 6709700
$ws.Range("F8").Value = 6570200
$ws.Range("G8").Value = 6324400
$ws.Range("H8").Value = 6136000
$ws.Range("I8").Value = 5824900
$ws.Range("J8").Value = 5368000

$ws.Range("D9").Value = 2639200
$ws.Range("E9").Value = 2731200
$ws.Range("F9").Value = 2552600
$ws.Range("G9").Value = 3523100
$ws.Range("H9").Value = 1648900
$ws.Range("I9").Value = 1676300
$ws.Range("J9").Value = 1556000

$ws.Range("D10").Value = 4181600
$ws.Range("E10").Value = 3978500
$ws.Range("F10").Value = 4017700
$ws.Range("G10").Value = 2801300
$ws.Range("H10").Value = 4487100
$ws.Range("I10").Value = 4148600
$ws.Range("J10").Value = 3812000

$ws.Range("D14").Value = 58100
$ws.Range("F14").Value = 54100
$ws.Range("I14").Value = 9700

$ws.Range("D15").Value = 1230600
$ws.Range("E15").Value = 1256400
$ws.Range("F15").Value = 1198800
$ws.Range("G15").Value = 1150900
$ws.Range("H15").Value = 1128100
$ws.Range("I15").Value = 1045800
$ws.Range("J15").Value = 984600

$ws.Range("D17").Value = 5340900
$ws.Range("E17").Value = 5359600
$ws.Range("F17").Value = 5123800
$ws.Range("G17").Value = 4786600
$ws.Range("H17").Value = 4568900
$ws.Range("I17").Value = 4380500
$ws.Range("J17").Value = 4131500

$ws.Range("D18").Value = 1479900
$ws.Range("E18").Value = 1350100
$ws.Range("F18").Value = 1446400
$ws.Range("G18").Value = 1537800
$ws.Range("H18").Value = 1567100
$ws.Range("I18").Value = 1444400
$ws.Range("J18").Value = 1236500

$ws.Range("D20").Value = 34600
$ws.Range("E20").Value = 52200
$ws.Range("F20").Value = -29200
$ws.Range("G20").Value = -344400
$ws.Range("H20").Value = -78100
$ws.Range("I20").Value = -168700
$ws.Range("J20").Value = -317800

$ws.Range("D21").Value = 2798400
$ws.Range("E21").Value = 2674000
$ws.Range("F21").Value = 2574100
$ws.Range("G21").Value = 2323000
$ws.Range("H21").Value = 2536200
$ws.Range("I21").Value = 2261600
$ws.Range("J21").Value = "NA"

$ws.Range("D22").Value = 350800
$ws.Range("E22").Value = 361100
$ws.Range("F22").Value = 385400
$ws.Range("G22").Value = 244400

$ws.Range("D23").Value = 1163700
$ws.Range("E23").Value = 1041200
$ws.Range("F23").Value = 1031800
$ws.Range("G23").Value = 949100
$ws.Range("H23").Value = 1489000
$ws.Range("I23").Value = 1275800
$ws.Range("J23").Value = 918700

$ws.Range("D24").Value = 292200
$ws.Range("E24").Value = 233100
$ws.Range("F24").Value = 214500
$ws.Range("G24").Value = 245000
$ws.Range("H24").Value = 302300
$ws.Range("I24").Value = 298500
$ws.Range("J24").Value = 239100

$ws.Range("D26").Value = 871500
$ws.Range("E26").Value = 808100
$ws.Range("F26").Value = 817300
$ws.Range("G26").Value = 704100
$ws.Range("H26").Value = 1186600
$ws.Range("I26").Value = 977300
$ws.Range("J26").Value = 679600

$ws.Range("D27").Value = 863000
$ws.Range("E27").Value = 808400
$ws.Range("F27").Value = 849400
$ws.Range("G27").Value = 701200
$ws.Range("H27").Value = 1172000
$ws.Range("I27").Value = 962300
$ws.Range("J27").Value = 623800

$ws.Range("E29").Value = -61900
$ws.Range("F29").Value = -87300
$ws.Range("G29").Value = 88900
$ws.Range("H29").Value = 57500
$ws.Range("I29").Value = -505800
$ws.Range("J29").Value = 27800

$ws.Range("D32").Value = -34600
$ws.Range("E32").Value = -52200
$ws.Range("F32").Value = 29200
$ws.Range("G32").Value = 344400
$ws.Range("H32").Value = 78100
$ws.Range("I32").Value = 168700
$ws.Range("J32").Value = 317800

$ws.Range("D33").Value = 863000
$ws.Range("E33").Value = 746500
$ws.Range("F33").Value = 762100
$ws.Range("G33").Value = 790100
$ws.Range("H33").Value = 1229500
$ws.Range("I33").Value = 456500
$ws.Range("J33").Value = 651700

$ws.Range("D35").Value = 863000
$ws.Range("E35").Value = 746500
$ws.Range("F35").Value = 762100
$ws.Range("G35").Value = 790100
$ws.Range("H35").Value = 1229500
$ws.Range("I35").Value = 456500
$ws.Range("J35").Value = 651700

$ws.Range("D41").Value = 471000
$ws.Range("E41").Value = 284400
$ws.Range("F41").Value = 515300
$ws.Range("G41").Value = 945700
$ws.Range("H41").Value = 471400
$ws.Range("I41").Value = 339000
$ws.Range("J41").Value = 28500

$ws.Range("D42").Value = 694400
$ws.Range("E42").Value = 117600
$ws.Range("F42").Value = 767500
$ws.Range("G42").Value = 153100
$ws.Range("H42").Value = 450700
$ws.Range("I42").Value = 62100

$ws.Range("D43").Value = 847700
$ws.Range("E43").Value = 676500
$ws.Range("F43").Value = 860400
$ws.Range("G43").Value = 1366700
$ws.Range("H43").Value = 1328000
$ws.Range("I43").Value = 602500
$ws.Range("J43").Value = 16300

$ws.Range("D44").Value = 153900
$ws.Range("E44").Value = 220700
$ws.Range("F44").Value = 223500
$ws.Range("G44").Value = 231300
$ws.Range("H44").Value = 261700
$ws.Range("I44").Value = 132500
$ws.Range("J44").Value = 4500

$ws.Range("D45").Value = 81900
$ws.Range("E45").Value = 102100
$ws.Range("F45").Value = 82100
$ws.Range("G45").Value = 448200
$ws.Range("H45").Value = 219700
$ws.Range("I45").Value = 287400
$ws.Range("J45").Value = 8500

$ws.Range("D46").Value = 2248900
$ws.Range("E46").Value = 1401300
$ws.Range("F46").Value = 2448900
$ws.Range("G46").Value = 2172200
$ws.Range("H46").Value = 1674600
$ws.Range("I46").Value = 1422400
$ws.Range("J46").Value = 59100

$ws.Range("D47").Value = 312600
$ws.Range("E47").Value = 992400
$ws.Range("F47").Value = 1116600
$ws.Range("G47").Value = 835200
$ws.Range("H47").Value = 366100
$ws.Range("I47").Value = 174700
$ws.Range("J47").Value = 4800

$ws.Range("D48").Value = 4057400
$ws.Range("E48").Value = 4206900
$ws.Range("F48").Value = 4666600
$ws.Range("G48").Value = 13826400
$ws.Range("H48").Value = 8358100
$ws.Range("I48").Value = 4185400
$ws.Range("J48").Value = 126900

$ws.Range("D49").Value = 1750600
$ws.Range("E49").Value = 1675700
$ws.Range("F49").Value = 1679600
$ws.Range("G49").Value = 4555600
$ws.Range("H49").Value = 2289300
$ws.Range("I49").Value = 1151900
$ws.Range("J49").Value = 41300

$ws.Range("D52").Value = 116900
$ws.Range("E52").Value = 108500
$ws.Range("F52").Value = 150400
$ws.Range("G52").Value = 460100
$ws.Range("H52").Value = 137100
$ws.Range("I52").Value = 94300

$ws.Range("D54").Value = 8486500
$ws.Range("E54").Value = 8384800
$ws.Range("F54").Value = 10062000
$ws.Range("G54").Value = 9229300
$ws.Range("H54").Value = 7364600
$ws.Range("I54").Value = 7006700
$ws.Range("J54").Value = 236000

$ws.Range("D57").Value = 745600
$ws.Range("E57").Value = 654300
$ws.Range("F57").Value = 917300
$ws.Range("G57").Value = 881800
$ws.Range("H57").Value = 655700
$ws.Range("I57").Value = 395300
$ws.Range("J57").Value = 25500

$ws.Range("D58").Value = 992900
$ws.Range("E58").Value = 727000
$ws.Range("F58").Value = 827000
$ws.Range("G58").Value = 987400
$ws.Range("H58").Value = 645600
$ws.Range("I58").Value = 423000
$ws.Range("J58").Value = 17800

$ws.Range("D59").Value = 674200
$ws.Range("E59").Value = 582900
$ws.Range("F59").Value = 687500
$ws.Range("G59").Value = 1441000
$ws.Range("H59").Value = 1280100
$ws.Range("I59").Value = 961000
$ws.Range("J59").Value = 29200

$ws.Range("D60").Value = 2412700
$ws.Range("E60").Value = 1964200
$ws.Range("F60").Value = 2431800
$ws.Range("G60").Value = 2146500
$ws.Range("H60").Value = 1563100
$ws.Range("I60").Value = 1746400
$ws.Range("J60").Value = 54900

$ws.Range("D61").Value = 3682100
$ws.Range("E61").Value = 3651500
$ws.Range("F61").Value = 4499400
$ws.Range("G61").Value = 3827700
$ws.Range("H61").Value = 2967000
$ws.Range("I61").Value = 3148300
$ws.Range("J61").Value = 116400

$ws.Range("D62").Value = 478900
$ws.Range("E62").Value = 552300
$ws.Range("F62").Value = 537900
$ws.Range("G62").Value = 656500
$ws.Range("H62").Value = 567500
$ws.Range("I62").Value = 323100
$ws.Range("J62").Value = 8500

$ws.Range("D66").Value = 6636500
$ws.Range("E66").Value = 6240600
$ws.Range("F66").Value = 7596300
$ws.Range("G66").Value = 6629300
$ws.Range("H66").Value = 5016200
$ws.Range("I66").Value = 5301500
$ws.Range("J66").Value = 182400

$ws.Range("D72").Value = 2189700
$ws.Range("E72").Value = 2149700
$ws.Range("F72").Value = 2839400
$ws.Range("G72").Value = 8438100
$ws.Range("H72").Value = 5582100
$ws.Range("I72").Value = 2284900
$ws.Range("J72").Value = 81500

$ws.Range("D76").Value = 1849900
$ws.Range("E76").Value = 2144200
$ws.Range("F76").Value = 2465800
$ws.Range("G76").Value = 2600000
$ws.Range("H76").Value = 2348400
$ws.Range("I76").Value = 1705100
$ws.Range("J76").Value = 53600

$ws.Range("D81").Value = 863000
$ws.Range("E81").Value = 746500
$ws.Range("F81").Value = 762100
$ws.Range("G81").Value = 790100
$ws.Range("H81").Value = 1229500
$ws.Range("I81").Value = 456500
$ws.Range("J81").Value = 651700

$ws.Range("D83").Value = 1282200
$ws.Range("E83").Value = 1270100
$ws.Range("F83").Value = 1155300
$ws.Range("G83").Value = 1128100
$ws.Range("H83").Value = 1045800
$ws.Range("I83").Value = 984600
$ws.Range("J83").Value = "NA"

$ws.Range("D89").Value = 2010700
$ws.Range("E89").Value = 2219000
$ws.Range("F89").Value = 2445900
$ws.Range("G89").Value = 2454400
$ws.Range("H89").Value = 2034700
$ws.Range("I89").Value = 1748900
$ws.Range("J89").Value = 55700

$ws.Range("D91").Value = -855300
$ws.Range("E91").Value = -1180700
$ws.Range("F91").Value = -2276400
$ws.Range("G91").Value = -1034000
$ws.Range("H91").Value = -1229500
$ws.Range("I91").Value = -985800
$ws.Range("J91").Value = -29500

$ws.Range("D94").Value = -882500
$ws.Range("E94").Value = -2238500
$ws.Range("F94").Value = -1617100
$ws.Range("G94").Value = -1488700
$ws.Range("H94").Value = -1437900
$ws.Range("I94").Value = -1189000
$ws.Range("J94").Value = "NA"

$ws.Range("D96").Value = -813200
$ws.Range("E96").Value = -782100
$ws.Range("F96").Value = -768800
$ws.Range("G96").Value = -611500
$ws.Range("H96").Value = -456200
$ws.Range("I96").Value = -535700
$ws.Range("J96").Value = -15000

$ws.Range("D100").Value = -1278800
$ws.Range("E100").Value = -425000
$ws.Range("F100").Value = -511500
$ws.Range("G100").Value = -849200
$ws.Range("H100").Value = -1160300
$ws.Range("I100").Value = -86700
$ws.Range("J100").Value = "NA"

$ws.Range("D101").Value = -80400
$ws.Range("E101").Value = 11700
$ws.Range("F101").Value = 157000
$ws.Range("G101").Value = 16000
$ws.Range("H101").Value = -15200
$ws.Range("I101").Value = 9100
$ws.Range("J101").Value = "NA"

$ws.Range("D102").Value = -230900
$ws.Range("E102").Value = -432800
$ws.Range("F102").Value = 474300
$ws.Range("G102").Value = 132400
$ws.Range("H102").Value = -578700
$ws.Range("I102").Value = 482300
$ws.Range("J102").Value = -24700
